$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix trailing whitespace typo in N4 ("2019-06-08  " -> "2019-06-08")
$ws.Range("N4").Value = "2019-06-08"

# 2. Fix typo in F5 ("...Urban Traill" -> "...Urban Trail")
$ws.Range("F5").Value = "Impression d'affiches, de flyer et de dépliants pour l'évènement Urban Trail"

# 3. Add new column V with header "anomalies" (extends used range to A1:V5)
$ws.Range("V1").Value = "anomalies"

# 4. Column width tweaks (closest achievable values through the character-width model)
#    Column F: 61.59 -> 61.18
$ws.Range("F1").ColumnWidth = 60.3
#    Column V (new): 9.77
$ws.Range("V1").ColumnWidth = 9.0
